$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the content of a Range with a fragment of WordprocessingML
# body content, using Range.InsertXML wrapped in a minimal OPC package. A
# fresh Range (via $d.Range(start, end)) is used as the target because the
# Range object returned directly by Find.Execute behaves like an insertion
# point for InsertXML rather than replacing its own span.
# ---------------------------------------------------------------------------
function Set-RangeXml($start, $end, $bodyXml) {
    $target = $d.Range($start, $end)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($pkg)
}

$rsquo = [char]0x2019
$mdash = [char]0x2014

# ---------------------------------------------------------------------------
# 1) "AdministratorAccess-Amplify'" run: drop the Consolas/shd "code" look
#    and give it Hyperlink character style (color auto, no underline)
#    instead.
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("AdministratorAccess-Amplify", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$end = $rng.End + 1   # include the trailing closing quote
$body = '<w:p><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="auto"/><w:u w:val="none"/></w:rPr>' +
        '<w:t>AdministratorAccess-Amplify' + $rsquo + '</w:t></w:r></w:p>'
Set-RangeXml $start $end $body

# ---------------------------------------------------------------------------
# 2) Merge the three runs "amplify-demo" + "-user" + "-iam" into a single
#    run, keeping only the Hyperlink character style (no explicit
#    color/underline overrides).
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("amplify-demo-user-iam", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$end = $rng.End
$body = '<w:p><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>amplify-demo-user-iam</w:t></w:r></w:p>'
Set-RangeXml $start $end $body

# ---------------------------------------------------------------------------
# 3) Replace the "Creates a profile7" bullet and the trailing empty bullet
#    with the full Angular/Amplify walk-through.
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("Creates a profile7", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Locate the paragraph containing the found text, and the one right after it
# (the empty trailing bullet before the section break) by scanning the
# paragraph collection for Range bounds that contain the found hit.
$startPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Start -le $rng.Start -and $pp.Range.End -ge $rng.End) {
        $startPara = $i
        break
    }
}
$p1 = $d.Paragraphs.Item($startPara)
$p2 = $d.Paragraphs.Item($startPara + 1)
$start = $p1.Range.Start
$end = $p2.Range.End

$body = ""
$body += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Install Amplify for Angular: </w:t></w:r></w:p>'
$body += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>npm install ' + $mdash + ' save aws-amplify</w:t></w:r></w:p>'
$body += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>npm install aws-amplify-angular</w:t></w:r></w:p>'
$body += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Initialize Amplify: amplify init</w:t></w:r></w:p>'
$body += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Does </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>lots</w:t></w:r><w:r><w:t xml:space="preserve"> of stuff so project is initialized and connected to the cloud</w:t></w:r></w:p>'
$body += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Add auth</w:t></w:r></w:p>'
$body += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Default config. Social provider needs Google and Facebook secrets</w:t></w:r></w:p>'
$body += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Use email to sign in</w:t></w:r></w:p>'
$body += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Advanced settings: Email and name to sign up</w:t></w:r></w:p>'
$body += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Email verification with redirect</w:t></w:r></w:p>'
$body += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>?????</w:t></w:r></w:p>'
$body += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Amplify push</w:t></w:r></w:p>'

Set-RangeXml $start $end $body

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
